$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.911.53'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '1.628.21'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.515'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.44%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.20'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.57%  '
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0607'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0882'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.27%  '
$ws.Range('D12').Value = '1.860.35'
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('D13').Value = '1.635.75'
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('E14').Value = '  -0.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.554'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.79'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.18%  '
$ws.Range('D17').Value = '27.924.46'
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '228.16'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.86%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.62'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('D20').Value = '0.0₃0717'
$ws.Range('E20').Value = '  -0.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.34'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.97'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.83%  '
$ws.Range('E24').Value = '  +1.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.25'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.47%  '
$ws.Range('E27').Value = '  -0.40%  '
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.41'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.88%  '
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.41'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('D33').Value = '1.418.46'
$ws.Range('E33').Value = '  +1.35%  '
$ws.Range('E34').Value = '  +0.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.62'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.58%  '
$ws.Range('E36').Value = '  -2.49%  '
$ws.Range('E37').Value = '  -0.61%  '
$ws.Range('E38').Value = '  -0.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.555'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.853'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.23%  '
$ws.Range('E41').Value = '  -2.58%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.82'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.83%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '65.60'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.21%  '
$ws.Range('E44').Value = '  -1.29%  '
$ws.Range('D45').Value = '1.769.74'
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('E46').Value = '  -3.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '88.74'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.74%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0102'
$ws.Range('E48').Value = '  -0.97%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.101'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.49%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0503'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.39%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.60'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.49%  '
